$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-03-01 03:18:27"
$ws.Range("H2").Value = "93%"
$ws.Range("N2").Value = "-1.0 °C 2:46 TU"
$ws.Range("O2").Value = "-0.5 °C"
$ws.Range("E3").Value = "2026-03-01 03:18:30"
$ws.Range("L3").Value = "16.6 km/h - 96º 2:46 TU"
$ws.Range("O3").Value = "-3.5 °C"
$ws.Range("E4").Value = "2026-03-01 03:18:32"
$ws.Range("J4").Value = "1025.6 hPa"
$ws.Range("M4").Value = "8.7 °C 2:34 TU"
$ws.Range("E5").Value = "2026-03-01 03:18:35"
$ws.Range("N5").Value = "-4.0 °C 2:58 TU"
$ws.Range("E6").Value = "2026-03-01 03:18:37"
$ws.Range("H6").Value = "83%"
$ws.Range("E7").Value = "2026-03-01 03:18:40"
$ws.Range("J7").Value = "1025.6 hPa"
$ws.Range("L7").Value = "10.8 km/h - 249º 2:58 TU"
$ws.Range("N7").Value = "13.3 °C 2:57 TU"
$ws.Range("E8").Value = "2026-03-01 03:18:42"
$ws.Range("M8").Value = "9.6 °C 2:35 TU"
$ws.Range("E9").Value = "2026-03-01 03:18:45"
$ws.Range("H9").Value = "62%"
$ws.Range("M9").Value = "12.0 °C 2:37 TU"
$ws.Range("O9").Value = "11.6 °C"
$ws.Range("E10").Value = "2026-03-01 03:18:47"
$ws.Range("O10").Value = "6.8 °C"
$ws.Range("E11").Value = "2026-03-01 03:18:49"
$ws.Range("N11").Value = "6.1 °C 2:30 TU"
$ws.Range("E12").Value = "2026-03-01 03:18:52"
$ws.Range("E13").Value = "2026-03-01 03:18:54"
$ws.Range("J13").Value = "1026.4 hPa"
$ws.Range("N13").Value = "4.3 °C 2:31 TU"
$ws.Range("E14").Value = "2026-03-01 03:18:56"
$ws.Range("N14").Value = "10.7 °C 2:43 TU"
$ws.Range("O14").Value = "11.3 °C"
$ws.Range("E15").Value = "2026-03-01 03:18:59"
$ws.Range("O15").Value = "8.1 °C"
$ws.Range("E16").Value = "2026-03-01 03:19:01"
$ws.Range("H16").Value = "86%"
$ws.Range("N16").Value = "-5.2 °C 2:58 TU"
$ws.Range("O16").Value = "-4.5 °C"
$ws.Range("E17").Value = "2026-03-01 03:19:04"
$ws.Range("N17").Value = "1.2 °C 2:59 TU"
$ws.Range("O17").Value = "1.5 °C"
$ws.Range("E18").Value = "2026-03-01 03:19:06"
$ws.Range("M18").Value = "7.7 °C 2:59 TU"
$ws.Range("O18").Value = "6.8 °C"
$ws.Range("E19").Value = "2026-03-01 03:19:09"
$ws.Range("N19").Value = "6.0 °C 2:50 TU"
$ws.Range("E20").Value = "2026-03-01 03:19:11"
$ws.Range("E21").Value = "2026-03-01 03:19:14"
$ws.Range("H21").Value = "87%"
$ws.Range("J21").Value = "1025.4 hPa"
$ws.Range("O21").Value = "6.8 °C"
$ws.Range("E22").Value = "2026-03-01 03:19:16"
$ws.Range("N22").Value = "-5.7 °C 2:58 TU"
$ws.Range("O22").Value = "-4.9 °C"
$ws.Range("E23").Value = "2026-03-01 03:19:18"
$ws.Range("N23").Value = "-3.7 °C 2:42 TU"
$ws.Range("O23").Value = "-3.3 °C"
$ws.Range("E24").Value = "2026-03-01 03:19:21"
$ws.Range("O24").Value = "4.3 °C"
$ws.Range("E25").Value = "2026-03-01 03:19:23"
$ws.Range("N25").Value = "-2.1 °C 2:59 TU"
$ws.Range("E26").Value = "2026-03-01 03:19:26"
$ws.Range("J26").Value = "1026.0 hPa"
$ws.Range("E27").Value = "2026-03-01 03:19:28"
$ws.Range("N27").Value = "-1.7 °C 2:48 TU"
$ws.Range("E28").Value = "2026-03-01 03:19:30"
$ws.Range("J28").Value = "1025.6 hPa"
$ws.Range("N28").Value = "8.4 °C 2:30 TU"
$ws.Range("E29").Value = "2026-03-01 03:19:33"
$ws.Range("H29").Value = "79%"
$ws.Range("N29").Value = "9.1 °C 2:35 TU"
$ws.Range("E30").Value = "2026-03-01 03:19:35"
$ws.Range("H30").Value = "84%"
$ws.Range("J30").Value = "1025.6 hPa"
$ws.Range("M30").Value = "11.5 °C 2:49 TU"
$ws.Range("O30").Value = "9.8 °C"
$ws.Range("E31").Value = "2026-03-01 03:19:38"
$ws.Range("J31").Value = "1024.5 hPa"
$ws.Range("N31").Value = "10.8 °C 2:56 TU"
$ws.Range("E32").Value = "2026-03-01 03:19:40"
$ws.Range("O32").Value = "1.6 °C"
$ws.Range("E33").Value = "2026-03-01 03:19:43"
$ws.Range("J33").Value = "1025.8 hPa"
$ws.Range("O33").Value = "4.0 °C"
$ws.Range("E34").Value = "2026-03-01 03:19:45"
$ws.Range("L34").Value = "7.6 km/h - 146º 2:48 TU"
$ws.Range("N34").Value = "-0.2 °C 2:49 TU"
$ws.Range("E35").Value = "2026-03-01 03:19:48"
$ws.Range("E36").Value = "2026-03-01 03:19:50"
$ws.Range("H36").Value = "80%"
$ws.Range("J36").Value = "1025.4 hPa"
$ws.Range("M36").Value = "10.3 °C 2:58 TU"
$ws.Range("O36").Value = "9.2 °C"
$ws.Range("E37").Value = "2026-03-01 03:19:52"
$ws.Range("L37").Value = "4.3 km/h - 15º 2:46 TU"
$ws.Range("N37").Value = "6.1 °C 2:43 TU"
$ws.Range("O37").Value = "6.4 °C"
$ws.Range("E38").Value = "2026-03-01 03:19:55"
$ws.Range("M38").Value = "9.4 °C 2:59 TU"
$ws.Range("O38").Value = "8.6 °C"
$ws.Range("E39").Value = "2026-03-01 03:19:57"
$ws.Range("E40").Value = "2026-03-01 03:19:59"
$ws.Range("H40").Value = "86%"
$ws.Range("N40").Value = "6.9 °C 2:48 TU"
$ws.Range("E41").Value = "2026-03-01 03:20:02"
$ws.Range("J41").Value = "1025.6 hPa"
$ws.Range("E42").Value = "2026-03-01 03:20:04"
$ws.Range("H42").Value = "79%"
$ws.Range("N42").Value = "7.3 °C 2:57 TU"
$ws.Range("O42").Value = "9.8 °C"
$ws.Range("E43").Value = "2026-03-01 03:20:06"
$ws.Range("H43").Value = "100%"
$ws.Range("N43").Value = "8.4 °C 2:59 TU"
$ws.Range("E44").Value = "2026-03-01 03:20:09"
$ws.Range("N44").Value = "-2.9 °C 2:55 TU"
$ws.Range("O44").Value = "-2.3 °C"
$ws.Range("E45").Value = "2026-03-01 03:20:11"
$ws.Range("N45").Value = "3.5 °C 2:59 TU"
$ws.Range("E46").Value = "2026-03-01 03:20:14"
$ws.Range("J46").Value = "1026.4 hPa"
$ws.Range("O46").Value = "7.4 °C"
